$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting of the adjacent G1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
